# Apply updated simulation results for case with 380 kV
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.3093457149082894
$ws.Range("C2").Value = 0.03562457795956675
$ws.Range("D2").Value = 0.07819256098835581
$ws.Range("E2").Value = 0.1566043186524908
$ws.Range("G2").Value = 0.002466008856346732
$ws.Range("K2").Value = 0.2680515022319412
$ws.Range("M2").Value = 0.2262230494872668
$ws.Range("O2").Value = 4.035369729367062
$ws.Range("B3").Value = 0.2789671693585376
$ws.Range("C3").Value = 0.03230477108915863
$ws.Range("D3").Value = 0.07097488807653463
$ws.Range("E3").Value = 0.1452944042285438
$ws.Range("G3").Value = 0.002469085992851357
$ws.Range("K3").Value = 0.2369903664306747
$ws.Range("M3").Value = 0.2050936524939928
$ws.Range("O3").Value = 4.018850526631297
$ws.Range("B4").Value = 0.2604118425794013
$ws.Range("C4").Value = 0.03024915261369188
$ws.Range("D4").Value = 0.06657856315327138
$ws.Range("E4").Value = 0.138449844219231
$ws.Range("G4").Value = 0.002471075522927404
$ws.Range("K4").Value = 0.2179566633124352
$ws.Range("M4").Value = 0.1922217994308397
$ws.Range("O4").Value = 4.010756527435376
$ws.Range("B5").Value = 0.252875065509528
$ws.Range("C5").Value = 0.02940714939460065
$ws.Range("D5").Value = 0.0647958939550648
$ws.Range("E5").Value = 0.1356854586358409
$ws.Range("G5").Value = 0.002471911536106263
$ws.Range("K5").Value = 0.2102100237026718
$ws.Range("M5").Value = 0.1870018235852058
$ws.Range("O5").Value = 4.007972706740162
$ws.Range("B6").Value = 0.2516250859922593
$ws.Range("C6").Value = 0.0292670747887982
$ws.Range("D6").Value = 0.06450041813538121
$ws.Range("E6").Value = 0.1352279264808374
$ws.Range("G6").Value = 0.002472051883485229
$ws.Range("K6").Value = 0.2089242949954979
$ws.Range("M6").Value = 0.1861365805796282
$ws.Range("O6").Value = 4.0075415139749
$ws.Range("B7").Value = 0.2603100987868174
$ws.Range("C7").Value = 0.03023781454533037
$ws.Range("D7").Value = 0.06655448557083332
$ws.Range("E7").Value = 0.1384124625630463
$ws.Range("G7").Value = 0.002471086695144962
$ws.Range("K7").Value = 0.2178521496753518
$ws.Range("M7").Value = 0.1921512983408249
$ws.Range("O7").Value = 4.010716901345774
$ws.Range("B8").Value = 0.29885113618991
$ws.Range("C8").Value = 0.03448349483669233
$ws.Range("D8").Value = 0.07569655410294729
$ws.Range("E8").Value = 0.1526837916720538
$ws.Range("G8").Value = 0.002467049112982659
$ws.Range("K8").Value = 0.2573338556152862
$ws.Range("M8").Value = 0.2189164259123615
$ws.Range("O8").Value = 4.029248260209414
$ws.Range("B9").Value = 0.375195794096328
$ws.Range("C9").Value = 0.04267221649608643
$ws.Range("D9").Value = 0.0939065922106721
$ws.Range("E9").Value = 0.1814747396847096
$ws.Range("G9").Value = 0.002459922529415093
$ws.Range("K9").Value = 0.3350536922977199
$ws.Range("M9").Value = 0.2722196279748204
$ws.Range("O9").Value = 4.0818811031981
$ws.Range("B10").Value = 0.4317516054355224
$ws.Range("C10").Value = 0.04860520482964148
$ws.Range("D10").Value = 0.1074617912405955
$ws.Range("E10").Value = 0.2031386717640444
$ws.Range("G10").Value = 0.002455163873537943
$ws.Range("K10").Value = 0.3923346087650259
$ws.Range("M10").Value = 0.3118975162484503
$ws.Range("O10").Value = 4.130542758941147
$ws.Range("B11").Value = 0.4575814752125211
$ws.Range("C11").Value = 0.05128629710905841
$ws.Range("D11").Value = 0.113667639140246
$ws.Range("E11").Value = 0.2131097200105714
$ws.Range("G11").Value = 0.002453101595500599
$ws.Range("K11").Value = 0.4184326158104739
$ws.Range("M11").Value = 0.3300642113666399
$ws.Range("O11").Value = 4.154863443649447
$ws.Range("B12").Value = 0.4673771684758492
$ws.Range("C12").Value = 0.05229898566336999
$ws.Range("D12").Value = 0.1160233606720169
$ws.Range("E12").Value = 0.2169025018174722
$ws.Range("G12").Value = 0.002452335316736831
$ws.Range("K12").Value = 0.4283209893332014
$ws.Range("M12").Value = 0.3369605407348075
$ws.Range("O12").Value = 4.164388054864105
$ws.Range("B13").Value = 0.4652668502839106
$ws.Range("C13").Value = 0.05208100018144535
$ws.Range("D13").Value = 0.1155157600143895
$ws.Range("E13").Value = 0.2160848995356588
$ws.Range("G13").Value = 0.002452499697736858
$ws.Range("K13").Value = 0.4261911044500835
$ws.Range("M13").Value = 0.3354745342086076
$ws.Range("O13").Value = 4.162322741540322
$ws.Range("B14").Value = 0.4583870838869188
$ws.Range("C14").Value = 0.05136966345871485
$ws.Range("D14").Value = 0.113861331393565
$ws.Range("E14").Value = 0.2134214129005159
$ws.Range("G14").Value = 0.002453038259943749
$ws.Range("K14").Value = 0.4192460267627496
$ws.Range("M14").Value = 0.3306312349358791
$ws.Range("O14").Value = 4.155640723562612
$ws.Range("B15").Value = 0.4541749073872268
$ws.Range("C15").Value = 0.05093361235304883
$ws.Range("D15").Value = 0.1128486891635845
$ws.Range("E15").Value = 0.2117921677204748
$ws.Range("G15").Value = 0.002453370052022451
$ws.Range("K15").Value = 0.4149926922226541
$ws.Range("M15").Value = 0.3276667927100192
$ws.Range("O15").Value = 4.151588827188874
$ws.Range("B16").Value = 0.4300656020182316
$ws.Range("C16").Value = 0.04842962853845734
$ws.Range("D16").Value = 0.1070570213293252
$ws.Range("E16").Value = 0.2024893987017506
$ws.Range("G16").Value = 0.002455300703930401
$ws.Range("K16").Value = 0.3906298415180629
$ws.Range("M16").Value = 0.3107126491750094
$ws.Range("O16").Value = 4.128997366338893
$ws.Range("B17").Value = 0.4153013694619005
$ws.Range("C17").Value = 0.04688893155993412
$ws.Range("D17").Value = 0.1035141560527961
$ws.Range("E17").Value = 0.1968123641559956
$ws.Range("G17").Value = 0.002456511286895131
$ws.Range("K17").Value = 0.3756942721654468
$ws.Range("M17").Value = 0.3003419232528444
$ws.Range("O17").Value = 4.115698247161646
$ws.Range("B18").Value = 0.4068190207260898
$ws.Range("C18").Value = 0.04600108381674772
$ws.Range("D18").Value = 0.1014801123437366
$ws.Range("E18").Value = 0.1935579873388917
$ws.Range("G18").Value = 0.002457217230065597
$ws.Range("K18").Value = 0.3671075649137947
$ws.Range("M18").Value = 0.2943879607739461
$ws.Range("O18").Value = 4.10825448966699
$ws.Range("B19").Value = 0.4039487073418115
$ws.Range("C19").Value = 0.04570018583463309
$ws.Range("D19").Value = 0.1007920582090236
$ws.Range("E19").Value = 0.1924579754318643
$ws.Range("G19").Value = 0.002457457909657145
$ws.Range("K19").Value = 0.3642009191144098
$ws.Range("M19").Value = 0.2923739363622673
$ws.Range("O19").Value = 4.10576943927353
$ws.Range("B20").Value = 0.4168720504736996
$ws.Range("C20").Value = 0.04705311532933365
$ws.Range("D20").Value = 0.103890915104401
$ws.Range("E20").Value = 0.1974155637081338
$ws.Range("G20").Value = 0.002456381420017251
$ws.Range("K20").Value = 0.3772837934736515
$ws.Range("M20").Value = 0.3014447645680249
$ws.Range("O20").Value = 4.117092682231402
$ws.Range("B21").Value = 0.4604074468655313
$ws.Range("C21").Value = 0.05157867051808296
$ws.Range("D21").Value = 0.1143471222941912
$ws.Range("E21").Value = 0.2142032809078174
$ws.Range("G21").Value = 0.002452879673782206
$ws.Range("K21").Value = 0.4212858135012709
$ws.Range("M21").Value = 0.3320533667153001
$ws.Range("O21").Value = 4.157594841628736
$ws.Range("B22").Value = 0.4889446509419031
$ws.Range("C22").Value = 0.05452132042762514
$ws.Range("D22").Value = 0.1212141001016107
$ws.Range("E22").Value = 0.2252740239236459
$ws.Range("G22").Value = 0.002450676503277674
$ws.Range("K22").Value = 0.4500763630598499
$ws.Range("M22").Value = 0.352157003704221
$ws.Range("O22").Value = 4.18590117153127
$ws.Range("B23").Value = 0.4737061711063859
$ws.Range("C23").Value = 0.05295215597658398
$ws.Range("D23").Value = 0.1175460170261289
$ws.Range("E23").Value = 0.2193562090041468
$ws.Range("G23").Value = 0.002451844584293406
$ws.Range("K23").Value = 0.4347073854054599
$ws.Range("M23").Value = 0.3414181864194887
$ws.Range("O23").Value = 4.170625306470527
$ws.Range("B24").Value = 0.4161619274896111
$ws.Range("C24").Value = 0.04697889432667068
$ws.Range("D24").Value = 0.1037205736291611
$ws.Range("E24").Value = 0.1971428278549752
$ws.Range("G24").Value = 0.002456440101604827
$ws.Range("K24").Value = 0.3765651710309896
$ws.Range("M24").Value = 0.3009461441549846
$ws.Range("O24").Value = 4.116461628691781
$ws.Range("B25").Value = 0.3544607188811995
$ws.Range("C25").Value = 0.04047159469232042
$ws.Range("D25").Value = 0.08894963158211056
$ws.Range("E25").Value = 0.1735977364016037
$ws.Range("G25").Value = 0.00246176628814171
$ws.Range("K25").Value = 0.313996756537307
$ws.Range("M25").Value = 0.2577103892198949
$ws.Range("O25").Value = 4.065892092262374
